# Préparation release concertation 220 ballot
# Updates Version / Date / Base Definition metadata and the
# Extension.value[x] Reference(...) string to pin the 2.2.0-ballot release.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
$meta.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# Base Definition: pin the FHIR core version
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s) column (K, row 6): pin the referenced profile version
$elements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-careplan-projet-personnalise|2.2.0-ballot)`n"

# Column K widened by the source tool's auto-fit after the longer text above
# (86.23828125 -> 95.74609375 "characters"); reproduce via the nearest width
# the host's pixel-grid column model can express.
$elements.Columns.Item(11).ColumnWidth = 94.75
